$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.631.46"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.510.24"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'586.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'183.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.49%  "
$ws.Range("D7").Value = "3.499.07"
$ws.Range("E7").Value = "  -2.65%  "
$ws.Range("D8").Value = "'0.611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.197"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.56%  "
$ws.Range("D11").Value = "'0.643"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").Value = "'53.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").Value = "'0.0000306"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "'9.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").Value = "4.064.05"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").Value = "'19.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").Value = "69.543.18"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "3.478.97"
$ws.Range("E18").Value = "  -3.37%  "
$ws.Range("D19").Value = "'12.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "'532.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.55%  "
$ws.Range("E22").Value = "  -3.58%  "
$ws.Range("D23").Value = "'18.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.18%  "
$ws.Range("D24").Value = "'4.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.16%  "
$ws.Range("D25").Value = "'4.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").Value = "'95.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'11.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'2.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").Value = "'9.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("D30").Value = "'32.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "'7.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.77%  "
$ws.Range("D32").Value = "'12.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("D33").Value = "'63.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("D35").Value = "'546.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.76%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'3.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.35%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "'0.407"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'38.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D40").Value = "0.0₃0759"
$ws.Range("E40").Value = "  -6.82%  "
$ws.Range("E41").Value = "  -2.11%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.348.60"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("D44").Value = "'3.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.73%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.29%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").Value = "'0.0437"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.134"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'8.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.22%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'138.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
